$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 2 ("Biometrics/Physical Unclonable Functions?")
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)

# 1) "TextBox 11" (Shapes.Item(6)) - change the word "prevent" to "inhibit" in
#    the last bullet, keeping the run/formatting intact.
$sh6 = $s2.Shapes.Item(6)
$tr6 = $sh6.TextFrame.TextRange
$para5 = $tr6.Paragraphs(5)
$fullPara = $tr6.Characters($para5.Start, $para5.Length)
$fullPara.Text = "Current techniques for removing noise impose large entropy losses and inhibit use in authentication"

# 2) "TextBox 14" (Shapes.Item(9)) - turn "Tuyls et al. 2006" into the
#    abbreviated citation "[Tuyls+06]" with a superscript "+".
$sh9 = $s2.Shapes.Item(9)
$tr9 = $sh9.TextFrame.TextRange

# "Tuyls" currently carries a stray spell-check marker (err="1"); drop it by
# deleting that run and retyping it together with the opening bracket, which
# inherits the clean formatting of the remaining " et al. 2006" text.
$tr9.Characters(1, 5).Delete()
$tr9.InsertBefore("[Tuyls")

# Replace " et al. " with "+".
$tr9.Characters(7, 8).Delete()
$tr9.Characters(7, 4).InsertBefore("+")

# Trim "2006" down to "06".
$tr9.Characters(8, 2).Delete()

# Close the citation with "]".
$tr9.InsertAfter("]")

# Re-assert the "[" text so it becomes its own run, separate from "Tuyls".
$tr9.Characters(1, 1).Text = "["

# Make the "+" a superscript, matching baseline="30000" in the XML.
$tr9.Characters(7, 1).Font.Superscript = 1

# The textbox auto-fits its width to the (now shorter) text; pin it to match
# the resized box from the authored edit (cx="1118365" EMU == ~88.0602pt).
$sh9.Width = 88.06027

# ---------------------------------------------------------------------------
# Slide 8 ("Computational Fuzzy Extractor")
# ---------------------------------------------------------------------------
$s8 = $p.Slides.Item(8)

# "Content Placeholder 3" (Shapes.Item(2)) - insert "secret key " right
# before the italic "w" in "Need encryption algorithm that allows decryption
# from close w1".
$sh2 = $s8.Shapes.Item(2)
$tr2 = $sh2.TextFrame.TextRange
$para3 = $tr2.Paragraphs(3)
$lead = $tr2.Characters($para3.Start, 60)
$lead.InsertAfter("secret key ")

# Force the inserted text into its own run (rather than merging into the
# preceding run) by re-asserting its (already matching) font name.
$newRun = $tr2.Characters($para3.Start + 60, 11)
$newRun.Font.Name = "Calibri"
